$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the target cells keep their text (string) representation,
# matching the workbook's original inlineStr cell type, instead of
# being auto-converted to numbers/percentages by Excel.
$cellUpdates = @{
    "D2" = "321.79"
    "E2" = "-2.76%"
    "D3" = "42.80"
    "E3" = "-5.76%"
    "D4" = "5.203"
    "E4" = "-6.79%"
    "D5" = "0.08192"
    "D6" = "4.304"
    "E6" = "-3.28%"
    "D7" = "1.814"
    "E7" = "-14.36%"
    "D8" = "0.9332"
    "E8" = "-4.52%"
    "D9" = "0.1110"
    "E9" = "-8.63%"
    "D10" = "0.1870"
    "E10" = "-2.44%"
    "D11" = "0.09430"
    "E11" = "-3.99%"
    "D12" = "0.04688"
    "E12" = "1.01%"
    "D13" = "7.422"
    "E13" = "-28.09%"
    "D14" = "0.1057"
    "E14" = "-0.01%"
    "D15" = "0.001288"
    "E15" = "-0.01%"
    "D16" = "0.005819"
    "E16" = "-1.61%"
    "D17" = "3.354"
    "E17" = "-0.64%"
    "D18" = "2.511"
    "E18" = "-2.56%"
    "E19" = "1.07%"
    "D20" = "0.1390"
    "E20" = "-0.19%"
    "E21" = "-8.57%"
    "D22" = "0.04145"
    "E22" = "-0.66%"
    "D23" = "0.001248"
    "E23" = "-3.63%"
    "D24" = "0.004350"
    "E24" = "-5.84%"
    "D25" = "0.0001200"
    "E25" = "-7.88%"
    "E26" = "-20.51%"
    "D38" = "0.02722"
    "E38" = "0.77%"
    "D39" = "0.05574"
    "E39" = "-2.84%"
    "D40" = "0.007984"
    "E40" = "1.22%"
    "D41" = "0.1398"
    "E41" = "-2.32%"
    "D42" = "0.006549"
    "E42" = "-12.89%"
    "D43" = "0.002093"
    "E43" = "-0.37%"
    "D44" = "0.008270"
    "E44" = "-2.86%"
    "D45" = "0.3492"
    "E45" = "3.63%"
    "D46" = "0.00006958"
    "E46" = "-2.39%"
    "E47" = "-0.20%"
    "D48" = "0.003510"
    "E48" = "-0.51%"
    "D49" = "0.003531"
    "E49" = "0.71%"
    "E50" = "-0.20%"
    "E51" = "-0.20%"
}

foreach ($addr in $cellUpdates.Keys) {
    $range = $ws.Range($addr)
    $range.NumberFormat = "@"
    $range.Value = $cellUpdates[$addr]
}
